$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("ALC")
$ws2 = $wb.Worksheets.Item("ARM")
$ws3 = $wb.Worksheets.Item("BSM")
$ws4 = $wb.Worksheets.Item("CRP")
$ws5 = $wb.Worksheets.Item("CUL")
$ws6 = $wb.Worksheets.Item("GSM")
$ws7 = $wb.Worksheets.Item("LTW")
$ws8 = $wb.Worksheets.Item("WVR")

# --- Sheet ALC ---
# row 15
$ws1.Range("H15").Value = 2061.92
$ws1.Range("I15").Value = 2061.92
$ws1.Range("K15").Value = 6185.76
$ws1.Range("M15").Value = -6016.76

# row 17
$ws1.Range("H17").Value = 1058334.2
$ws1.Range("J17").Value = 1137619.8
$ws1.Range("L17").Value = 3412859.4
$ws1.Range("N17").Value = -3413195.4

# row 40
$ws1.Range("H40").Value = 2271
$ws1.Range("I40").Value = 1999.5
$ws1.Range("J40").Value = 2379.6
$ws1.Range("K40").Value = 1999.5
$ws1.Range("L40").Value = 2379.6
$ws1.Range("M40").Value = -1824.5
$ws1.Range("N40").Value = -2729.6

# row 62
$ws1.Range("H62").Value = 22533.555
$ws1.Range("I62").Value = 4933.8887
$ws1.Range("J62").Value = 40133.223
$ws1.Range("K62").Value = 4933.8887
$ws1.Range("L62").Value = 40133.223
$ws1.Range("M62").Value = -4309.8887
$ws1.Range("N62").Value = -41381.223

# row 65
$ws1.Range("H65").Value = 22533.555
$ws1.Range("I65").Value = 4933.8887
$ws1.Range("J65").Value = 40133.223
$ws1.Range("K65").Value = 24669.4435
$ws1.Range("L65").Value = 200666.115
$ws1.Range("M65").Value = -21549.4435
$ws1.Range("N65").Value = -206906.115

# row 132
$ws1.Range("H132").Value = 5046.877
$ws1.Range("I132").Value = 5125.8237
$ws1.Range("K132").Value = 15377.4711
$ws1.Range("M132").Value = -12847.4711

# row 137
$ws1.Range("H137").Value = 5080.302
$ws1.Range("I137").Value = 6760.122
$ws1.Range("K137").Value = 20280.366
$ws1.Range("M137").Value = -17730.366

# row 138
$ws1.Range("H138").Value = 414222.4
$ws1.Range("J138").Value = 3873.6667
$ws1.Range("L138").Value = 11621.0001
$ws1.Range("N138").Value = -21901.0001

# row 141
$ws1.Range("H141").Value = 3687.3215
$ws1.Range("I141").Value = 3453.5186
$ws1.Range("J141").Value = 10000
$ws1.Range("K141").Value = 10360.5558
$ws1.Range("L141").Value = 30000
$ws1.Range("M141").Value = -5180.5558
$ws1.Range("N141").Value = -40360

# --- Sheet ARM ---
# row 32
$ws2.Range("H32").Value = 2396.9
$ws2.Range("I32").Value = 2396.9
$ws2.Range("J32").Value = 0
$ws2.Range("K32").Value = 2396.9
$ws2.Range("L32").Value = 0
$ws2.Range("M32").Value = -2109.9
$ws2.Range("N32").ClearContents()

# row 61
$ws2.Range("H61").Value = 3640.1897
$ws2.Range("I61").Value = 3292.5305
$ws2.Range("J61").Value = 5533
$ws2.Range("K61").Value = 3292.5305
$ws2.Range("L61").Value = 5533
$ws2.Range("M61").Value = -3080.5305
$ws2.Range("N61").Value = -5957

# row 102
$ws2.Range("H102").Value = 12428.714
$ws2.Range("I102").Value = 15063.0625
$ws2.Range("K102").Value = 15063.0625
$ws2.Range("M102").Value = -13441.0625

# row 122
$ws2.Range("H122").Value = 2005981.4
$ws2.Range("I122").Value = 8875.429
$ws2.Range("J122").Value = 2937964.2
$ws2.Range("K122").Value = 26626.287
$ws2.Range("L122").Value = 8813892.600000001
$ws2.Range("M122").Value = -24176.287
$ws2.Range("N122").Value = -8818792.600000001

# row 132
$ws2.Range("H132").Value = 6395.75
$ws2.Range("I132").Value = 6590.227
$ws2.Range("J132").Value = 4256.5
$ws2.Range("K132").Value = 19770.681
$ws2.Range("L132").Value = 12769.5
$ws2.Range("M132").Value = -17240.681
$ws2.Range("N132").Value = -17829.5

# row 136
$ws2.Range("H136").Value = 3640.1897
$ws2.Range("I136").Value = 3292.5305
$ws2.Range("J136").Value = 5533
$ws2.Range("K136").Value = 9877.591499999999
$ws2.Range("L136").Value = 16599
$ws2.Range("M136").Value = -7327.591499999999
$ws2.Range("N136").Value = -21699

# --- Sheet BSM ---
# row 22
$ws3.Range("H22").Value = 0
$ws3.Range("J22").Value = 0
$ws3.Range("L22").Value = 0
$ws3.Range("N22").ClearContents()

# row 95
$ws3.Range("H95").Value = 41849.75
$ws3.Range("J95").Value = 41849.75
$ws3.Range("L95").Value = 41849.75
$ws3.Range("N95").Value = -47341.75

# row 105
$ws3.Range("H105").Value = 62863.234
$ws3.Range("I105").Value = 92561.37
$ws3.Range("J105").Value = 8416.666999999999
$ws3.Range("K105").Value = 92561.37
$ws3.Range("L105").Value = 8416.666999999999
$ws3.Range("M105").Value = -90814.37
$ws3.Range("N105").Value = -11910.667

# row 134
$ws3.Range("H134").Value = 11132.481
$ws3.Range("I134").Value = 12107.458
$ws3.Range("J134").Value = 3332.6667
$ws3.Range("K134").Value = 36322.374
$ws3.Range("L134").Value = 9998.000100000001
$ws3.Range("M134").Value = -33787.374
$ws3.Range("N134").Value = -15068.0001

# --- Sheet CRP ---
# row 7
$ws4.Range("H7").Value = 6644.7744
$ws4.Range("I7").Value = 11240.777
$ws4.Range("J7").Value = 281.07693
$ws4.Range("K7").Value = 11240.777
$ws4.Range("L7").Value = 281.07693
$ws4.Range("M7").Value = -11127.777
$ws4.Range("N7").Value = -507.07693

# row 16
$ws4.Range("H16").Value = 987.1667
$ws4.Range("I16").Value = 987.1667
$ws4.Range("K16").Value = 987.1667
$ws4.Range("M16").Value = -700.1667

# row 22
$ws4.Range("H22").Value = 1176
$ws4.Range("I22").Value = 1257.4
$ws4.Range("J22").Value = 1108.1666
$ws4.Range("K22").Value = 1257.4
$ws4.Range("L22").Value = 1108.1666
$ws4.Range("M22").Value = -907.4000000000001
$ws4.Range("N22").Value = -1808.1666

# row 31
$ws4.Range("H31").Value = 3153.3958
$ws4.Range("I31").Value = 2437.513
$ws4.Range("J31").Value = 6255.5557
$ws4.Range("K31").Value = 2437.513
$ws4.Range("L31").Value = 6255.5557
$ws4.Range("M31").Value = -2142.513
$ws4.Range("N31").Value = -6845.5557

# row 34
$ws4.Range("H34").Value = 3153.3958
$ws4.Range("I34").Value = 2437.513
$ws4.Range("J34").Value = 6255.5557
$ws4.Range("K34").Value = 2437.513
$ws4.Range("L34").Value = 6255.5557
$ws4.Range("M34").Value = -2235.513
$ws4.Range("N34").Value = -6659.5557

# row 58
$ws4.Range("H58").Value = 2710.7
$ws4.Range("I58").Value = 2315.5715
$ws4.Range("J58").Value = 3632.6667
$ws4.Range("K58").Value = 2315.5715
$ws4.Range("L58").Value = 3632.6667
$ws4.Range("M58").Value = -2112.5715
$ws4.Range("N58").Value = -4038.6667

# row 82
$ws4.Range("H82").Value = 45999.8
$ws4.Range("I82").Value = 38749.75
$ws4.Range("K82").Value = 38749.75
$ws4.Range("M82").Value = -38388.75

# row 85
$ws4.Range("H85").Value = 45999.8
$ws4.Range("I85").Value = 38749.75
$ws4.Range("K85").Value = 38749.75
$ws4.Range("M85").Value = -37501.75

# row 88
$ws4.Range("H88").Value = 39993.668
$ws4.Range("J88").Value = 37992.4
$ws4.Range("L88").Value = 37992.4
$ws4.Range("N88").Value = -38804.4

# row 91
$ws4.Range("H91").Value = 39993.668
$ws4.Range("J91").Value = 37992.4
$ws4.Range("L91").Value = 37992.4
$ws4.Range("N91").Value = -40800.4

# row 105
$ws4.Range("H105").Value = 107031.3
$ws4.Range("I105").Value = 141922.6
$ws4.Range("J105").Value = 2357.4
$ws4.Range("K105").Value = 141922.6
$ws4.Range("L105").Value = 2357.4
$ws4.Range("M105").Value = -140175.6
$ws4.Range("N105").Value = -5851.4

# row 107
$ws4.Range("H107").Value = 11158.909
$ws4.Range("I107").Value = 18998.666
$ws4.Range("K107").Value = 18998.666
$ws4.Range("M107").Value = -17078.666

# row 113
$ws4.Range("H113").Value = 987.1667
$ws4.Range("I113").Value = 987.1667
$ws4.Range("K113").Value = 987.1667
$ws4.Range("M113").Value = 1182.8333

# row 122
$ws4.Range("H122").Value = 4112.6045
$ws4.Range("I122").Value = 5142
$ws4.Range("J122").Value = 1737.0769
$ws4.Range("K122").Value = 15426
$ws4.Range("L122").Value = 5211.2307
$ws4.Range("M122").Value = -12976
$ws4.Range("N122").Value = -10111.2307

# row 134
$ws4.Range("H134").Value = 3161.1538
$ws4.Range("I134").Value = 2623
$ws4.Range("K134").Value = 7869
$ws4.Range("M134").Value = -5334

# row 136
$ws4.Range("H136").Value = 2710.7
$ws4.Range("I136").Value = 2315.5715
$ws4.Range("J136").Value = 3632.6667
$ws4.Range("K136").Value = 6946.7145
$ws4.Range("L136").Value = 10898.0001
$ws4.Range("M136").Value = -4396.7145
$ws4.Range("N136").Value = -15998.0001

# --- Sheet CUL ---
# row 2
$ws5.Range("H2").Value = 138.90909
$ws5.Range("J2").Value = 216.66667
$ws5.Range("L2").Value = 1300.00002
$ws5.Range("N2").Value = -1526.00002

# row 129
$ws5.Range("H129").Value = 22808200
$ws5.Range("I129").Value = 932.44446
$ws5.Range("J129").Value = 43334744
$ws5.Range("K129").Value = 2797.33338
$ws5.Range("L129").Value = 130004232
$ws5.Range("M129").Value = 2202.66662
$ws5.Range("N129").Value = -130014232

# row 137
$ws5.Range("H137").Value = 4082
$ws5.Range("I137").Value = 3450
$ws5.Range("J137").Value = 4924.6665
$ws5.Range("K137").Value = 10350
$ws5.Range("L137").Value = 14773.9995
$ws5.Range("M137").Value = -5250
$ws5.Range("N137").Value = -24973.9995

# --- Sheet GSM ---
# row 97
$ws6.Range("H97").Value = 5046.121
$ws6.Range("I97").Value = 5362.5
$ws6.Range("J97").Value = 3871
$ws6.Range("K97").Value = 5362.5
$ws6.Range("L97").Value = 3871
$ws6.Range("M97").Value = -4866.5
$ws6.Range("N97").Value = -4863

# row 102
$ws6.Range("H102").Value = 6511.971
$ws6.Range("I102").Value = 7456.4814
$ws6.Range("K102").Value = 7456.4814
$ws6.Range("M102").Value = -5834.4814

# row 123
$ws6.Range("H123").Value = 0
$ws6.Range("J123").Value = 0
$ws6.Range("L123").Value = 0
$ws6.Range("N123").ClearContents()

# row 132
$ws6.Range("H132").Value = 4570.7383
$ws6.Range("J132").Value = 7655.1665
$ws6.Range("L132").Value = 22965.4995
$ws6.Range("N132").Value = -28025.4995

# --- Sheet LTW ---
# row 16
$ws7.Range("H16").Value = 2491.6
$ws7.Range("I16").Value = 2491.6
$ws7.Range("K16").Value = 2491.6
$ws7.Range("M16").Value = -2321.6

# row 40
$ws7.Range("H40").Value = 25758.9
$ws7.Range("I40").Value = 38999
$ws7.Range("J40").Value = 9576.556
$ws7.Range("K40").Value = 38999
$ws7.Range("L40").Value = 9576.556
$ws7.Range("M40").Value = -38863
$ws7.Range("N40").Value = -9848.556

# row 46
$ws7.Range("H46").Value = 3000
$ws7.Range("I46").Value = 2714.2856
$ws7.Range("J46").Value = 5000
$ws7.Range("K46").Value = 2714.2856
$ws7.Range("L46").Value = 5000
$ws7.Range("M46").Value = -2526.2856
$ws7.Range("N46").Value = -5376

# row 136
$ws7.Range("H136").Value = 6813.4736
$ws7.Range("I136").Value = 3141.4546
$ws7.Range("K136").Value = 9424.363799999999
$ws7.Range("M136").Value = -6874.363799999999

# --- Sheet WVR ---
# row 92
$ws8.Range("H92").Value = 110041000
$ws8.Range("J92").Value = 110041000
$ws8.Range("L92").Value = 110041000
$ws8.Range("N92").Value = -110045992

# row 132
$ws8.Range("H132").Value = 14650.956
$ws8.Range("I132").Value = 16942.916
$ws8.Range("J132").Value = 6399.9
$ws8.Range("K132").Value = 50828.74800000001
$ws8.Range("L132").Value = 19199.7
$ws8.Range("M132").Value = -48298.74800000001
$ws8.Range("N132").Value = -24259.7
